$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column G (old Urbanizado column),
# shifting old G -> H and old H -> I
$ws.Columns.Item(7).Insert()

# New header for inserted column G
$ws.Range("G1").Value = "Uso_solo_simplificado"

# Populate simplified land-use classification per row based on old Urbanizado/Transporte flags
$ws.Range("G2").Value = "Outros"
$ws.Range("G3").Value = "Urbanizado"
$ws.Range("G4").Value = "Outros"
$ws.Range("G5").Value = "Transporte"
$ws.Range("G6").Value = "Outros"
$ws.Range("G7").Value = "Outros"
$ws.Range("G8").Value = "Outros"
$ws.Range("G9").Value = "Outros"
$ws.Range("G10").Value = "Urbanizado"
$ws.Range("G11").Value = "Urbanizado"
$ws.Range("G12").Value = "Outros"
$ws.Range("G13").Value = "Urbanizado"
$ws.Range("G14").Value = "Outros"
$ws.Range("G15").Value = "Transporte"
$ws.Range("G16").Value = "Urbanizado"
$ws.Range("G17").Value = "Outros"
$ws.Range("G18").Value = "Outros"
$ws.Range("G19").Value = "Outros"
$ws.Range("G20").Value = "Outros"
$ws.Range("G21").Value = "Outros"
$ws.Range("G22").Value = "Transporte"
$ws.Range("G23").Value = "Outros"
$ws.Range("G24").Value = "Outros"
$ws.Range("G25").Value = "Outros"
$ws.Range("G26").Value = "Outros"
$ws.Range("G27").Value = "Urbanizado"
$ws.Range("G28").Value = "Outros"
$ws.Range("G29").Value = "Outros"
$ws.Range("G30").Value = "Urbanizado"
$ws.Range("G31").Value = "Outros"
$ws.Range("G32").Value = "Urbanizado"
$ws.Range("G33").Value = "Urbanizado"
$ws.Range("G34").Value = "Urbanizado"
$ws.Range("G35").Value = "Urbanizado"
$ws.Range("G36").Value = "Outros"
$ws.Range("G37").Value = "Outros"
$ws.Range("G38").Value = "Outros"
$ws.Range("G39").Value = "Outros"
$ws.Range("G40").Value = "Outros"
$ws.Range("G41").Value = "Urbanizado"
$ws.Range("G42").Value = "Urbanizado"
$ws.Range("G43").Value = "Urbanizado"
$ws.Range("G44").Value = "Outros"
$ws.Range("G45").Value = "Outros"
$ws.Range("G46").Value = "Outros"
$ws.Range("G47").Value = "Outros"
$ws.Range("G48").Value = "Outros"
$ws.Range("G49").Value = "Outros"
